$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 116
$ws_ALC.Range("H116").Value = 1836931.2
$ws_ALC.Range("I116").Value = 5954010
$ws_ALC.Range("K116").Value = 5954010
$ws_ALC.Range("M116").Value = -5950568

# ALC row 132
$ws_ALC.Range("H132").Value = 2598839.5
$ws_ALC.Range("I132").Value = 2756024
$ws_ALC.Range("J132").Value = 5298
$ws_ALC.Range("K132").Value = 8268072
$ws_ALC.Range("L132").Value = 15894
$ws_ALC.Range("M132").Value = -8265542
$ws_ALC.Range("N132").Value = -20954

# ARM row 61
$ws_ARM.Range("H61").Value = 1426.25
$ws_ARM.Range("I61").Value = 1235
$ws_ARM.Range("K61").Value = 1235
$ws_ARM.Range("M61").Value = -1023

# ARM row 112
$ws_ARM.Range("H112").Value = 21666.666
$ws_ARM.Range("J112").Value = 21666.666
$ws_ARM.Range("L112").Value = 21666.666
$ws_ARM.Range("N112").Value = -24620.666

# ARM row 114
$ws_ARM.Range("H114").Value = 0
$ws_ARM.Range("J114").Value = 0
$ws_ARM.Range("L114").Value = 0
$ws_ARM.Range("N114").ClearContents()

# ARM row 132
$ws_ARM.Range("H132").Value = 1556.2693
$ws_ARM.Range("I132").Value = 1587.0454
$ws_ARM.Range("J132").Value = 1387
$ws_ARM.Range("K132").Value = 4761.1362
$ws_ARM.Range("L132").Value = 4161
$ws_ARM.Range("M132").Value = -2231.1362
$ws_ARM.Range("N132").Value = -9221

# ARM row 136
$ws_ARM.Range("H136").Value = 1426.25
$ws_ARM.Range("I136").Value = 1235
$ws_ARM.Range("K136").Value = 3705
$ws_ARM.Range("M136").Value = -1155

# BSM row 20
$ws_BSM.Range("H20").Value = 22732074
$ws_BSM.Range("I20").Value = 34488436
$ws_BSM.Range("J20").Value = 3106.6
$ws_BSM.Range("K20").Value = 34488436
$ws_BSM.Range("L20").Value = 3106.6
$ws_BSM.Range("M20").Value = -34488189
$ws_BSM.Range("N20").Value = -3600.6

# BSM row 134
$ws_BSM.Range("H134").Value = 1204.3829
$ws_BSM.Range("I134").Value = 956.8889
$ws_BSM.Range("J134").Value = 1538.5
$ws_BSM.Range("K134").Value = 2870.6667
$ws_BSM.Range("L134").Value = 4615.5
$ws_BSM.Range("M134").Value = -335.6667000000002
$ws_BSM.Range("N134").Value = -9685.5

# CRP row 31
$ws_CRP.Range("H31").Value = 1828.4225
$ws_CRP.Range("I31").Value = 1129.6724
$ws_CRP.Range("J31").Value = 4945.923
$ws_CRP.Range("K31").Value = 1129.6724
$ws_CRP.Range("L31").Value = 4945.923
$ws_CRP.Range("M31").Value = -834.6723999999999
$ws_CRP.Range("N31").Value = -5535.923

# CRP row 34
$ws_CRP.Range("H34").Value = 1828.4225
$ws_CRP.Range("I34").Value = 1129.6724
$ws_CRP.Range("J34").Value = 4945.923
$ws_CRP.Range("K34").Value = 1129.6724
$ws_CRP.Range("L34").Value = 4945.923
$ws_CRP.Range("M34").Value = -927.6723999999999
$ws_CRP.Range("N34").Value = -5349.923

# CRP row 96
$ws_CRP.Range("H96").Value = 21333.818
$ws_CRP.Range("J96").Value = 21333.818
$ws_CRP.Range("L96").Value = 21333.818
$ws_CRP.Range("N96").Value = -26825.818

# CRP row 99
$ws_CRP.Range("H99").Value = 3000
$ws_CRP.Range("I99").Value = 2500
$ws_CRP.Range("J99").Value = 5000
$ws_CRP.Range("K99").Value = 2500
$ws_CRP.Range("L99").Value = 5000
$ws_CRP.Range("M99").Value = -1002
$ws_CRP.Range("N99").Value = -7996

# CRP row 126
$ws_CRP.Range("H126").Value = 3000
$ws_CRP.Range("I126").Value = 2500
$ws_CRP.Range("J126").Value = 5000
$ws_CRP.Range("K126").Value = 7500
$ws_CRP.Range("L126").Value = 15000
$ws_CRP.Range("M126").Value = -5030
$ws_CRP.Range("N126").Value = -19940

# CRP row 132
$ws_CRP.Range("H132").Value = 1552.7441
$ws_CRP.Range("I132").Value = 1273.1936
$ws_CRP.Range("J132").Value = 2274.9167
$ws_CRP.Range("K132").Value = 3819.5808
$ws_CRP.Range("L132").Value = 6824.750100000001
$ws_CRP.Range("M132").Value = -1289.5808
$ws_CRP.Range("N132").Value = -11884.7501

# CUL row 5
$ws_CUL.Range("H5").Value = 349.5
$ws_CUL.Range("I5").Value = 226.36667
$ws_CUL.Range("J5").Value = 811.25
$ws_CUL.Range("K5").Value = 679.10001
$ws_CUL.Range("L5").Value = 2433.75
$ws_CUL.Range("M5").Value = -567.10001
$ws_CUL.Range("N5").Value = -2657.75

# CUL row 76
$ws_CUL.Range("H76").Value = 6192.3076
$ws_CUL.Range("J76").Value = 7150
$ws_CUL.Range("L76").Value = 21450
$ws_CUL.Range("N76").Value = -22216

# CUL row 79
$ws_CUL.Range("H79").Value = 6192.3076
$ws_CUL.Range("J79").Value = 7150
$ws_CUL.Range("L79").Value = 21450
$ws_CUL.Range("N79").Value = -24102

# CUL row 95
$ws_CUL.Range("H95").Value = 8081.75
$ws_CUL.Range("J95").Value = 8081.75
$ws_CUL.Range("L95").Value = 24245.25
$ws_CUL.Range("N95").Value = -28363.25

# CUL row 113
$ws_CUL.Range("H113").Value = 4410167.5
$ws_CUL.Range("I113").Value = 568049.1
$ws_CUL.Range("J113").Value = 17857580
$ws_CUL.Range("K113").Value = 1704147.3
$ws_CUL.Range("L113").Value = 53572740
$ws_CUL.Range("M113").Value = -1701977.3
$ws_CUL.Range("N113").Value = -53577080

# CUL row 122
$ws_CUL.Range("H122").Value = 1163195.6
$ws_CUL.Range("I122").Value = 326.2857
$ws_CUL.Range("J122").Value = 6250749
$ws_CUL.Range("K122").Value = 2936.5713
$ws_CUL.Range("L122").Value = 56256741
$ws_CUL.Range("M122").Value = -486.5713000000001
$ws_CUL.Range("N122").Value = -56261641

# CUL row 131
$ws_CUL.Range("H131").Value = 846.53845
$ws_CUL.Range("J131").Value = 1039.963
$ws_CUL.Range("L131").Value = 3119.889
$ws_CUL.Range("N131").Value = -13199.889

# CUL row 135
$ws_CUL.Range("H135").Value = 349.5
$ws_CUL.Range("I135").Value = 226.36667
$ws_CUL.Range("J135").Value = 811.25
$ws_CUL.Range("K135").Value = 2037.30003
$ws_CUL.Range("L135").Value = 7301.25
$ws_CUL.Range("M135").Value = 497.6999700000001
$ws_CUL.Range("N135").Value = -12371.25

# GSM row 132
$ws_GSM.Range("H132").Value = 5127.033
$ws_GSM.Range("I132").Value = 5530.654
$ws_GSM.Range("J132").Value = 2503.5
$ws_GSM.Range("K132").Value = 16591.962
$ws_GSM.Range("L132").Value = 7510.5
$ws_GSM.Range("M132").Value = -14061.962
$ws_GSM.Range("N132").Value = -12570.5

# LTW row 7
$ws_LTW.Range("H7").Value = 4985.7144
$ws_LTW.Range("I7").Value = 4725
$ws_LTW.Range("J7").Value = 5333.3335
$ws_LTW.Range("K7").Value = 4725
$ws_LTW.Range("L7").Value = 5333.3335
$ws_LTW.Range("M7").Value = -4613
$ws_LTW.Range("N7").Value = -5557.3335

# LTW row 106
$ws_LTW.Range("H106").Value = 18623.334
$ws_LTW.Range("J106").Value = 18623.334
$ws_LTW.Range("L106").Value = 18623.334
$ws_LTW.Range("N106").Value = -21147.334

# LTW row 126
$ws_LTW.Range("H126").Value = 4985.7144
$ws_LTW.Range("I126").Value = 4725
$ws_LTW.Range("J126").Value = 5333.3335
$ws_LTW.Range("K126").Value = 14175
$ws_LTW.Range("L126").Value = 16000.0005
$ws_LTW.Range("M126").Value = -11705
$ws_LTW.Range("N126").Value = -20940.0005

# WVR row 56
$ws_WVR.Range("H56").Value = 11179.375
$ws_WVR.Range("I56").Value = 992.5
$ws_WVR.Range("J56").Value = 14575
$ws_WVR.Range("K56").Value = 992.5
$ws_WVR.Range("L56").Value = 14575
$ws_WVR.Range("M56").Value = -278.5
$ws_WVR.Range("N56").Value = -16003

# WVR row 105
$ws_WVR.Range("H105").Value = 32211.8
$ws_WVR.Range("J105").Value = 32211.8
$ws_WVR.Range("L105").Value = 32211.8
$ws_WVR.Range("N105").Value = -39199.8
